$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "massachusetts"
$ws.Range("H2").Value = "00-0001"
$ws.Range("I2").Value = "test"
$ws.Range("J2").Value = "major"
$ws.Range("K2").Value = 5
$ws.Range("Y2").Value = 5
